# Apply the updated crypto price / volume values described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new text, and whether it must be forced to
# plain text (column D holds numeric-looking strings such as "25.772.24"
# or "0.3830" whose exact formatting -- thousands separators, trailing
# zeros, etc. -- must be preserved rather than re-interpreted as a number).
$updates = @(
    @{ Cell = "D2"; Value = "25.772.24"; ForceText = $true },
    @{ Cell = "E2"; Value = "  -2.29%  "; ForceText = $false },
    @{ Cell = "D3"; Value = "1.752.26"; ForceText = $true },
    @{ Cell = "E3"; Value = "  -4.30%  "; ForceText = $false },
    @{ Cell = "E4"; Value = "  -0.10%  "; ForceText = $false },
    @{ Cell = "D5"; Value = "236.94"; ForceText = $true },
    @{ Cell = "E5"; Value = "  -6.47%  "; ForceText = $false },
    @{ Cell = "E6"; Value = "  -0.23%  "; ForceText = $false },
    @{ Cell = "D7"; Value = "0.5076"; ForceText = $true },
    @{ Cell = "E7"; Value = "  -3.27%  "; ForceText = $false },
    @{ Cell = "D8"; Value = "41.55"; ForceText = $true },
    @{ Cell = "E8"; Value = "  -6.35%  "; ForceText = $false },
    @{ Cell = "D9"; Value = "0.2646"; ForceText = $true },
    @{ Cell = "E9"; Value = "  -4.37%  "; ForceText = $false },
    @{ Cell = "D10"; Value = "0.06165"; ForceText = $true },
    @{ Cell = "E10"; Value = "  -9.56%  "; ForceText = $false },
    @{ Cell = "D11"; Value = "1.756.27"; ForceText = $true },
    @{ Cell = "E11"; Value = "  -4.23%  "; ForceText = $false },
    @{ Cell = "D12"; Value = "15.74"; ForceText = $true },
    @{ Cell = "E12"; Value = "  -4.44%  "; ForceText = $false },
    @{ Cell = "D13"; Value = "0.06917"; ForceText = $true },
    @{ Cell = "E13"; Value = "  -2.67%  "; ForceText = $false },
    @{ Cell = "D14"; Value = "0.6047"; ForceText = $true },
    @{ Cell = "E14"; Value = "  -11.36%  "; ForceText = $false },
    @{ Cell = "D15"; Value = "4.505"; ForceText = $true },
    @{ Cell = "E15"; Value = "  -6.92%  "; ForceText = $false },
    @{ Cell = "D16"; Value = "77.19"; ForceText = $true },
    @{ Cell = "E16"; Value = "  -10.30%  "; ForceText = $false },
    @{ Cell = "E17"; Value = "  -0.16%  "; ForceText = $false },
    @{ Cell = "D18"; Value = "1.001"; ForceText = $true },
    @{ Cell = "E18"; Value = "  -0.13%  "; ForceText = $false },
    @{ Cell = "D19"; Value = "25.779.43"; ForceText = $true },
    @{ Cell = "E19"; Value = "  -2.36%  "; ForceText = $false },
    @{ Cell = "D20"; Value = "0.000006868"; ForceText = $true },
    @{ Cell = "E20"; Value = "  -6.10%  "; ForceText = $false },
    @{ Cell = "D21"; Value = "11.71"; ForceText = $true },
    @{ Cell = "E21"; Value = "  -11.04%  "; ForceText = $false },
    @{ Cell = "D22"; Value = "1.975.03"; ForceText = $true },
    @{ Cell = "E22"; Value = "  -5.29%  "; ForceText = $false },
    @{ Cell = "D23"; Value = "4.092"; ForceText = $true },
    @{ Cell = "E23"; Value = "  -8.67%  "; ForceText = $false },
    @{ Cell = "D24"; Value = "8.259"; ForceText = $true },
    @{ Cell = "E24"; Value = "  -7.72%  "; ForceText = $false },
    @{ Cell = "D25"; Value = "5.213"; ForceText = $true },
    @{ Cell = "D26"; Value = "137.67"; ForceText = $true },
    @{ Cell = "E26"; Value = "  -3.05%  "; ForceText = $false },
    @{ Cell = "D27"; Value = "1.471"; ForceText = $true },
    @{ Cell = "E27"; Value = "  -12.31%  "; ForceText = $false },
    @{ Cell = "E29"; Value = "  -8.65%  "; ForceText = $false },
    @{ Cell = "D30"; Value = "102.89"; ForceText = $true },
    @{ Cell = "E30"; Value = "  -5.39%  "; ForceText = $false },
    @{ Cell = "D31"; Value = "0.08212"; ForceText = $true },
    @{ Cell = "E31"; Value = "  -5.74%  "; ForceText = $false },
    @{ Cell = "E32"; Value = "  -8.55%  "; ForceText = $false },
    @{ Cell = "E33"; Value = "  -9.52%  "; ForceText = $false },
    @{ Cell = "D34"; Value = "0.04523"; ForceText = $true },
    @{ Cell = "E34"; Value = "  -3.53%  "; ForceText = $false },
    @{ Cell = "E35"; Value = "  -0.11%  "; ForceText = $false },
    @{ Cell = "D36"; Value = "2.658"; ForceText = $true },
    @{ Cell = "E36"; Value = "  -7.64%  "; ForceText = $false },
    @{ Cell = "D37"; Value = "0.9995"; ForceText = $true },
    @{ Cell = "E37"; Value = "  -9.07%  "; ForceText = $false },
    @{ Cell = "D38"; Value = "0.6089"; ForceText = $true },
    @{ Cell = "E38"; Value = "  -12.90%  "; ForceText = $false },
    @{ Cell = "D39"; Value = "2.695"; ForceText = $true },
    @{ Cell = "E39"; Value = "  -11.78%  "; ForceText = $false },
    @{ Cell = "D40"; Value = "1.951"; ForceText = $true },
    @{ Cell = "E40"; Value = "  -9.78%  "; ForceText = $false },
    @{ Cell = "D41"; Value = "0.01556"; ForceText = $true },
    @{ Cell = "E41"; Value = "  -4.70%  "; ForceText = $false },
    @{ Cell = "E42"; Value = "  -0.15%  "; ForceText = $false },
    @{ Cell = "D43"; Value = "103.63"; ForceText = $true },
    @{ Cell = "E43"; Value = "  -1.22%  "; ForceText = $false },
    @{ Cell = "D44"; Value = "0.3830"; ForceText = $true },
    @{ Cell = "E44"; Value = "  -13.49%  "; ForceText = $false },
    @{ Cell = "D45"; Value = "0.7403"; ForceText = $true },
    @{ Cell = "E45"; Value = "  -14.22%  "; ForceText = $false },
    @{ Cell = "D46"; Value = "4.927"; ForceText = $true },
    @{ Cell = "E46"; Value = "  -13.87%  "; ForceText = $false },
    @{ Cell = "D47"; Value = "0.05466"; ForceText = $true },
    @{ Cell = "E47"; Value = "  -1.90%  "; ForceText = $false },
    @{ Cell = "E48"; Value = "  -5.54%  "; ForceText = $false },
    @{ Cell = "D49"; Value = "6.020"; ForceText = $true },
    @{ Cell = "E49"; Value = "  -13.69%  "; ForceText = $false },
    @{ Cell = "D50"; Value = "7.706"; ForceText = $true },
    @{ Cell = "E50"; Value = "  -10.25%  "; ForceText = $false },
    @{ Cell = "D51"; Value = "29.99"; ForceText = $true },
    @{ Cell = "E51"; Value = "  -9.84%  "; ForceText = $false }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Leading apostrophe forces Excel to store the value as text instead
        # of parsing it as a number (which would drop trailing zeros / turn
        # multi-dot strings into errors, or round the value).
        $range.Value = "'" + $u.Value
        # Revert to the default style so no extra "@" text format is left
        # applied to the cell (keeps styling identical to the original).
        $range.Style = "Normal"
    } else {
        $range.Value = $u.Value
    }
}
